$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the sheet is active (mirrors tabSelected="1" / the user working in it)
$ws.Activate()

# Enter the "last report 29-03-25" figures for itopup (C2C) and Distributor
# GA & SAF Commission rows. C17 (Total income) and F43 (Profit/Loss) are
# driven by formulas and recalculate automatically.
$ws.Range("C13").Value = 2174
$ws.Range("C15").Value = 5400

# Reposition the view the way it was left after the edit: scrolled so row 34
# is at the top, with C16 as the active selected cell.
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C16").Select()
